$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.632.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "'1.590.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'211.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "'1.814.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'1.602.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "'64.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'26.636.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").Value = "'208.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'2.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'146.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "'0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("D33").Value = "'0.665"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +22.43%  "
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").Value = "'1.318.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").Value = "'0.789"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'63.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "'1.726.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").Value = "'90.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("D49").Value = "'0.0512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'0.0979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  -0.19%  "
